$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A92").Value = "$ 27.553 CLP 26-11-20"
